# Update column C ("Förändrad") date value from 2023-09-20 (45189) to
# 2023-09-21 (45190) for every data row (rows 2 through 135).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}
